$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (rich-text runs with identical formatting across runs) ---
$ws.Range("A8").Value = "Volume 30   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/27/2023  Through  12/3/2023"

# --- Type changes: text placeholder <-> numeric, using format+value paste so the
#     resulting style index matches a same-row donor cell exactly (avoids minting new xf). ---

# Row 22: D22 (text "N/A") -> number 1 ; E22 (text "***.*") -> number 0 (pct-style)
$ws.Range("C22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null   # xlPasteValues (type stays number)
$ws.Range("D22").Value = 1
$ws.Range("H22").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("H22").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = 0

# Row 23: C23 (number 1) -> text "N/A" placeholder
$ws.Range("D23").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4163) | Out-Null

# Row 27: D27 (text) -> number 1 ; E27 (text) -> number 100
$ws.Range("C27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("D27").Value = 1
$ws.Range("H27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("H27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = 100

# Row 28: D28 (text) -> number 1 ; E28 (text) -> number -100
$ws.Range("F28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("F28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("D28").Value = 1
$ws.Range("H28").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("H28").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4163) | Out-Null
$ws.Range("E28").Value = -100

# Row 29: D29 (text) -> number 1 ; E29 (text) -> number -100
$ws.Range("F29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("F29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("D29").Value = 1
$ws.Range("H29").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("H29").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = -100

# Row 30: D30 (number 1) -> text "N/A" placeholder ; E30 (number -100) -> text "***.*" placeholder
$ws.Range("C30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("C30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("M30").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("M30").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false

# --- Plain numeric value updates ---
$ws.Range("M15").Value = -61.111111111111
$ws.Range("N15").Value = -41.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 23.076923076923
$ws.Range("I16").Value = 161
$ws.Range("J16").Value = 189
$ws.Range("K16").Value = -14.814814814814
$ws.Range("L16").Value = -26.48401826484
$ws.Range("M16").Value = -28.125
$ws.Range("N16").Value = 62.626262626262
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 7.692307692307
$ws.Range("I17").Value = 229
$ws.Range("J17").Value = 252
$ws.Range("K17").Value = -9.126984126984
$ws.Range("L17").Value = -19.366197183098
$ws.Range("M17").Value = 36.309523809523
$ws.Range("N17").Value = 197.402597402597
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -38.461538461538
$ws.Range("I18").Value = 108
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = 9.090909090909
$ws.Range("L18").Value = -13.6
$ws.Range("M18").Value = 9.090909090909
$ws.Range("N18").Value = -8.474576271186
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 9.523809523809
$ws.Range("I19").Value = 331
$ws.Range("J19").Value = 275
$ws.Range("K19").Value = 20.363636363636
$ws.Range("L19").Value = 0.30303030303
$ws.Range("M19").Value = 32.4
$ws.Range("N19").Value = 401.515151515152
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 28.571428571428
$ws.Range("I20").Value = 133
$ws.Range("J20").Value = 127
$ws.Range("K20").Value = 4.724409448818
$ws.Range("L20").Value = 75
$ws.Range("M20").Value = 107.8125
$ws.Range("N20").Value = 47.777777777777
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 6.25
$ws.Range("F21").Value = 70
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = 4.477611940298
$ws.Range("I21").Value = 976
$ws.Range("J21").Value = 957
$ws.Range("K21").Value = 1.985370950888
$ws.Range("L21").Value = -7.047619047619
$ws.Range("M21").Value = 18.016928657799
$ws.Range("N21").Value = 111.255411255411
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 350
$ws.Range("I22").Value = 27
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 42.105263157894
$ws.Range("M22").Value = 92.857142857142
$ws.Range("L23").Value = -39.130434782608
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -26.666666666666
$ws.Range("F24").Value = 51
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = -28.169014084507
$ws.Range("I24").Value = 855
$ws.Range("J24").Value = 1217
$ws.Range("K24").Value = -29.74527526705
$ws.Range("L24").Value = -42.655935613682
$ws.Range("M24").Value = 55.454545454545
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = -24
$ws.Range("I25").Value = 345
$ws.Range("J25").Value = 372
$ws.Range("K25").Value = -7.258064516129
$ws.Range("L25").Value = -10.38961038961
$ws.Range("M25").Value = -23.333333333333
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 40
$ws.Range("J27").Value = 47
$ws.Range("K27").Value = -14.893617021276
$ws.Range("L27").Value = -21.56862745098
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 12
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = 9.090909090909
$ws.Range("L28").Value = -33.333333333333
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 8
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = -20
$ws.Range("L29").Value = -50
